# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-08-12 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-08-13 Wednesday", 2)

# Update every arithmetic-problem cell in the single 20-row x 5-column table.
# Row indices below refer to the rows that actually carry text
# (rows 1, 5, 9, 13, 17 of the 20-row table); columns are 1-5.
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; New="91÷6="},
    @{Row=1;  Col=2; New="36÷3="},
    @{Row=1;  Col=3; New="37÷8="},
    @{Row=1;  Col=4; New="87÷5="},
    @{Row=1;  Col=5; New="90÷5="},

    @{Row=5;  Col=1; New="27÷6="},
    @{Row=5;  Col=2; New="11÷3="},
    @{Row=5;  Col=3; New="65÷4="},
    @{Row=5;  Col=4; New="27÷6="},
    @{Row=5;  Col=5; New="47÷9="},

    @{Row=9;  Col=1; New="80÷6="},
    @{Row=9;  Col=2; New="21÷6="},
    @{Row=9;  Col=3; New="70÷5="},
    @{Row=9;  Col=4; New="47÷4="},
    @{Row=9;  Col=5; New="17÷4="},

    @{Row=13; Col=1; New="22÷6="},
    @{Row=13; Col=2; New="11÷8="},
    @{Row=13; Col=3; New="98÷7="},
    @{Row=13; Col=4; New="46÷6="},
    @{Row=13; Col=5; New="67÷7="},

    @{Row=17; Col=1; New="43÷9="},
    @{Row=17; Col=2; New="65÷5="},
    @{Row=17; Col=3; New="36÷4="},
    @{Row=17; Col=4; New="68÷8="},
    @{Row=17; Col=5; New="87÷5="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}
